$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: add new column F header ("c_source" moves from E to F) and
#    reorder B/C/D/E into c_office_chn, c_dy, c_office_trans, c_office_pinyin.
#    A1 already carries the bordered/bold header style (s="1"); clone that
#    format onto B1:F1 before overwriting the header captions.
# ---------------------------------------------------------------------------
$ws.Range("A2:E2").ClearContents()

$ws.Range("A1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$headers = @("c_office_id", "c_office_chn", "c_dy", "c_office_trans", "c_office_pinyin", "c_source")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---------------------------------------------------------------------------
# 2. Data rows: c_office_id, c_office_chn, c_dy, c_office_trans,
#    c_office_pinyin, c_source
# ---------------------------------------------------------------------------
$data = @()
$data += ,@(802012, "僧錄司", "19", "The Central Buddhist Registry", "seng lu si", "18417")
$data += ,@(802013, "僧錄司善世", "19", "Buddhist Patriarch of the Central Buddhist Registry", "seng lu si shan shi", "18417")
$data += ,@(802014, "僧錄司左善世", "19", "Left Buddhist Patriarch of the Central Buddhist Registry", "seng lu si zuo shan shi", "18417")
$data += ,@(802015, "僧錄司右善世", "19", "Right Buddhist Patriarch of the Central Buddhist Registry", "seng lu si you shan shi", "18417")
$data += ,@(802016, "僧錄司闡教", "19", "Supervisory Monk for Buddhist Practices of the Central Buddhist Registry", "seng lu si chan jiao", "18417")
$data += ,@(802017, "僧錄司左闡教", "19", "Left Supervisory Monk for Buddhist Practices of the Central Buddhist Registry", "seng lu si zuo chan jiao", "18417")
$data += ,@(802018, "僧錄司右闡教", "19", "Right Supervisory Monk for Buddhist Practices of the Central Buddhist Registry", "seng lu si you chan jiao", "18417")
$data += ,@(802019, "僧錄司講經", "19", "Lecturing Monk of the Central Buddhist Registry", "seng lu si jiang jing", "18417")
$data += ,@(802020, "僧錄司左講經", "19", "Left Lecturing Monk of the Central Buddhist Registry", "seng lu si zuo jiang jing", "18417")
$data += ,@(802021, "僧錄司右講經", "19", "Right Lecturing Monk of the Central Buddhist Registry", "seng lu si you jiang jing", "18417")
$data += ,@(802022, "僧錄司覺義", "19", "Buddhist Rectifier of the Central Buddhist Registry", "seng lu si jue yi", "18417")
$data += ,@(802023, "僧錄司左覺義", "19", "Left Buddhist Rectifier of the Central Buddhist Registry", "seng lu si zuo jue yi", "18417")
$data += ,@(802024, "僧錄司右覺義", "19", "Right Buddhist Rectifier of the Central Buddhist Registry", "seng lu si you jue yi", "18417")
$data += ,@(802025, "府學訓導", "20", "Assistant Instructor in a Prefectural Confucian School", "fu xue xun dao", "65006")

# Helper cell (far outside the used range) used to coerce numeric-looking
# text ("19", "20", "18417", "65006") into real text values instead of
# numbers, the same way Excel requires a Text-formatted cell for that.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

function Set-TextValue($cell, $text) {
    $helper.Value = $text
    $helper.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    Set-TextValue $ws.Cells.Item($row, 3) $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    Set-TextValue $ws.Cells.Item($row, 6) $record[5]
    $row++
}

$helper.EntireColumn.Delete()
$ws.Application.CutCopyMode = $false
